{"js": "// Add the \"OSM als Beispiel\" idea and a \"Referenzen\" note after the\n// existing \"Hier vielzahl an Datenquellen ...\" paragraph, each new\n// paragraph preceded by a blank line. The _GoBack bookmark (Word's\n// \"last edit position\" marker) moves along with the edit, from the\n// old final paragraph to the new \"OSM als Beispiel\" paragraph.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// The _GoBack bookmark currently sits at the end of the last paragraph;\n// drop it so it can be re-inserted at the new edit location below.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Blank paragraph, then the \"OSM als Beispiel\" paragraph.\nconst blank1 = lastParagraph.insertParagraph(\"\", \"After\");\nconst osmParagraph = blank1.insertParagraph(\"OSM als Beispiel\", \"After\");\n\n// Blank paragraph, then the \"Referenzen\" paragraph.\nconst blank2 = osmParagraph.insertParagraph(\"\", \"After\");\nblank2.insertParagraph(\"Referenzen: z.B. GISCUP (gewinner), etc\", \"After\");\n\n// Re-create the _GoBack bookmark over the newly typed \"OSM als Beispiel\" text.\nconst osmRange = osmParagraph.getRange(\"Content\");\nosmRange.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Add the \"OSM als Beispiel\" idea and a \"Referenzen\" note after the\n# existing \"Hier vielzahl an Datenquellen ...\" paragraph, each new\n# paragraph preceded by a blank line. The _GoBack bookmark (Word's\n# \"last edit position\" marker) moves along with the edit, from the\n# old final paragraph to the new \"OSM als Beispiel\" paragraph.\n\n$d = $word.ActiveDocument\n\n# The _GoBack bookmark currently sits at the end of the last paragraph;\n# drop it so it can be re-inserted at the new edit location below.\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n\n# Blank paragraph after the current last paragraph.\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n\n# Another blank paragraph, which we then fill with \"OSM als Beispiel\".\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$osmIndex = $d.Paragraphs.Count\n$d.Paragraphs.Item($osmIndex).Range.InsertAfter(\"OSM als Beispiel\")\n\n# Blank paragraph after that.\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n\n# Another blank paragraph, which we then fill with the \"Referenzen\" text.\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$refIndex = $d.Paragraphs.Count\n$d.Paragraphs.Item($refIndex).Range.InsertAfter(\"Referenzen: z.B. GISCUP (gewinner), etc\")\n\n# Re-create the _GoBack bookmark over the newly typed \"OSM als Beispiel\" text\n# (exclude the trailing paragraph mark from the bookmark range).\n$osmRange = $d.Paragraphs.Item($osmIndex).Range\n[void]$osmRange.MoveEnd(1, -1)\n$bookmarkRange = $d.Range($osmRange.Start, $osmRange.End)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange) | Out-Null\n"}
